# Apply updated regression estimates to censoring_imp sheet
# Cells are formula cells referencing closed external workbooks
# (e.g. =[1]decomposition_main_te_0_0!B5). Since the external sources
# are not available to recompute, we refresh each cell's cached value
# by rewriting it as a literal-string formula (="value") which keeps
# the cell a formula/text (t="str") cell with the updated cached result,
# matching the type and layout of the original cells as closely as possible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Formula = '="-234.6***"'
$ws.Range("C5").Formula = '="-191.8***"'
$ws.Range("D5").Formula = '="0.089"'
$ws.Range("E5").Formula = '="-74.5**"'
$ws.Range("B10").Formula = '="988.6"'
$ws.Range("D10").Formula = '="5.26"'
$ws.Range("E10").Formula = '="395.2"'
$ws.Range("F10").Formula = '="0.43"'
$ws.Range("B14").Formula = '="-189.9***"'
$ws.Range("D14").Formula = '="1.89"'
$ws.Range("E14").Formula = '="-13.8"'
$ws.Range("F14").Formula = '="0.0094"'
$ws.Range("G14").Formula = '="-0.073***"'
$ws.Range("G18").Formula = '="0.022"'
$ws.Range("B19").Formula = '="988.6"'
$ws.Range("D19").Formula = '="5.26"'
$ws.Range("E19").Formula = '="395.2"'
$ws.Range("F19").Formula = '="0.43"'
$ws.Range("B41").Formula = '="-262.6***"'
$ws.Range("C41").Formula = '="-170.7***"'
$ws.Range("D41").Formula = '="-0.56"'
$ws.Range("E41").Formula = '="-123.9***"'
$ws.Range("B42").Formula = '="(53.3)"'
$ws.Range("D42").Formula = '="(3.30)"'
$ws.Range("E42").Formula = '="(32.5)"'
$ws.Range("B43").Formula = '="-39.5"'
$ws.Range("C43").Formula = '="-30.5"'
$ws.Range("D43").Formula = '="-2.71"'
$ws.Range("E43").Formula = '="-10.4"'
$ws.Range("F43").Formula = '="-0.014"'
$ws.Range("G43").Formula = '="0.0060"'
$ws.Range("D44").Formula = '="(2.75)"'
$ws.Range("E44").Formula = '="(34.4)"'
$ws.Range("B47").Formula = '="0.017"'
$ws.Range("D47").Formula = '="0.003"'
$ws.Range("E47").Formula = '="0.009"'
$ws.Range("F47").Formula = '="0.016"'
$ws.Range("G47").Formula = '="0.042"'
